# Updated symbol list on Fri Feb  3 22:26:44 UTC 2023 with GitHub Actions
# Applies the latest coinranking.com snapshot values to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in columns D (Price) and E (Volume(1h)) hold numeric-looking / percent-looking
# text (e.g. "329.27", "0.35%"). These must stay plain text, matching the source
# workbook's inlineStr cells, so force Text number format before writing the value -
# otherwise Excel auto-converts them to numbers/percentages.
$textCells = @{
    'D2' = '329.27';
    'E2' = '0.35%';
    'D3' = '41.31';
    'E3' = '2.81%';
    'D4' = '5.717';
    'E4' = '-1.87%';
    'D5' = '0.08073';
    'E5' = '0.56%';
    'D6' = '2.048';
    'E6' = '4.45%';
    'D7' = '8.718';
    'E7' = '0.01%';
    'D8' = '4.514';
    'E8' = '-1.63%';
    'D9' = '2.921';
    'E9' = '-0.69%';
    'D10' = '0.9210';
    'E10' = '-2.56%';
    'D11' = '0.1241';
    'E11' = '-0.94%';
    'D12' = '0.1951';
    'E12' = '-0.50%';
    'D13' = '8.273';
    'E13' = '-6.97%';
    'D14' = '0.09365';
    'E14' = '1.60%';
    'D15' = '0.03674';
    'E15' = '2.68%';
    'D16' = '0.1052';
    'E16' = '9.22%';
    'D17' = '0.001297';
    'E17' = '-0.19%';
    'D18' = '0.006198';
    'E18' = '2.20%';
    'E19' = '0.40%';
    'D20' = '0.3482';
    'E20' = '-1.30%';
    'E21' = '0.82%';
    'D22' = '0.2651';
    'E22' = '9.68%';
    'D23' = '0.04435';
    'E23' = '0.74%';
    'D24' = '0.001257';
    'E24' = '-0.39%';
    'D25' = '0.004374';
    'E25' = '1.28%';
    'E26' = '8.42%';
    'D39' = '0.02831';
    'E39' = '16.67%';
    'D40' = '0.05475';
    'E40' = '3.61%';
    'D41' = '0.007591';
    'E41' = '1.98%';
    'D42' = '0.009942';
    'E42' = '16.91%';
    'D43' = '0.1417';
    'E43' = '0.16%';
    'D44' = '0.002119';
    'E44' = '0.50%';
    'D45' = '0.01188';
    'E45' = '8.76%';
    'D46' = '0.00006746';
    'E46' = '-2.11%';
    'E47' = '-0.36%';
    'D48' = '0.002990';
    'E48' = '-5.28%';
    'D49' = '0.002280';
    'E49' = '59.91%';
    'D50' = '0.00002101';
    'E50' = '-0.36%';
    'D51' = '0.0002001';
    'E51' = '-0.36%'
}

foreach ($addr in $textCells.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $textCells[$addr]
}

# Columns B (Coin) and C (Link) are plain (non-numeric) text, no special formatting needed.
$plainCells = @{
    'B48' = 'BOLO';
    'C48' = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo';
    'B49' = 'CoinbaseStockToken';
    'C49' = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
}

foreach ($addr in $plainCells.Keys) {
    $ws.Range($addr).Value = $plainCells[$addr]
}
